$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): update "想去人数" (column F) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1337
$wsExpo.Range("F3").Value = 2858

# Sheet "全部类型" (all types): same two events appear one row further down
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1337
$wsAll.Range("F4").Value = 2858
